# Insert a new event row at row 29 ("Un viaggio alla scoperta della musica"),
# pushing the existing rows 29-59 down to 30-60, and populate the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 29; this shifts rows 29:59 -> 30:60
# and Excel auto-extends the used range / dimension to A1:AA60.
$ws.Rows("29:29").Insert()

$newRow = $ws.Range("A29:AA29")

# Y29/Z29 hold comma-decimal coordinate text ("44,648090007085386") that Excel
# would otherwise auto-parse as a (huge) number; force Text format on just
# those two cells before writing so they are stored as literal text like every
# other text cell on the sheet.
$ws.Range("Y29:Z29").NumberFormat = "@"

$ws.Cells.Item(29, 1).Value2 = "Iniziative per bambini"
$ws.Cells.Item(29, 2).Value2 = "Modena"
$ws.Cells.Item(29, 3).Value2 = "piazza Matteotti, 17"
$ws.Cells.Item(29, 4).Value2 = "2022-05-25T09:39:41+00:00"
$ws.Cells.Item(29, 5).Value2 = "Spettacolo teatrale per bambini e bambine da 3 a 6 anni"
$ws.Cells.Item(29, 6).Value2 = "2022-05-25T09:43:08+00:00"
$ws.Cells.Item(29, 7).Value2 = ""
$ws.Cells.Item(29, 8).Value2 = "2022-05-29T08:00:00+00:00"
$ws.Cells.Item(29, 9).Value2 = "2022-05-29T10:00:00+00:00"
$ws.Cells.Item(29, 10).Value2 = "https://www.comune.modena.it/api/novita/eventi/2022/un-viaggio-alla-scoperta-della-musica/@@images/c4306e19-428c-42cb-bc32-698e6ea5dd39.jpeg"
$ws.Cells.Item(29, 11).Value2 = ""
$ws.Cells.Item(29, 12).Value2 = "2022-05-25T09:43:08+00:00"
$ws.Cells.Item(29, 13).Value2 = "Centro per l'infanzia Mo.Mo"
$ws.Cells.Item(29, 14).Value2 = " Ore 10.00"
$ws.Cells.Item(29, 15).Value2 = ""
$ws.Cells.Item(29, 16).Value2 = " Inizia gratuita su prenotazione"
$ws.Cells.Item(29, 17).Value2 = ""
$ws.Cells.Item(29, 18).Value2 = "059 235320"
$ws.Cells.Item(29, 19).Value2 = "Un viaggio alla scoperta della musica"
$ws.Cells.Item(29, 20).Value2 = ""
$ws.Cells.Item(29, 21).Value2 = ""
$ws.Cells.Item(29, 22).Value2 = $false
$ws.Cells.Item(29, 23).Value2 = 41123
$ws.Cells.Item(29, 24).Value2 = "https://www.comune.modena.it/novita/eventi/2022/un-viaggio-alla-scoperta-della-musica"
$ws.Cells.Item(29, 25).Value2 = "44,648090007085386"
$ws.Cells.Item(29, 26).Value2 = "10,925069996586558"
$ws.Cells.Item(29, 27).Value2 = "POINT (10.925069996586558 44.648090007085386)"

# Drop the transient Text-format style from Y29:Z29 now that the literal
# values are committed, so the new row matches the rest of the sheet's
# default (unstyled) data cells.
$newRow.ClearFormats()
